$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the checkpoint #1 date in cell B16 from "Friday, July 15th, 2022"
# to "Friday, July 22nd, 2022" (one week later)
$ws.Range("B16").Value = "Friday, July 22nd, 2022"

# Update the view: scroll the sheet so row 9 is the top-left visible row,
# change zoom from 145% to 115%, and move the selection to E22
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("E22").Select()
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Application.ActiveWindow.ScrollColumn = 1
